# Apply edits to the "Contract Details" sheet per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C (Excel's ColumnWidth<->stored XML width has a fixed +5/6
# offset in this engine, so subtract it to land on the exact stored width).
$ws.Columns.Item(3).ColumnWidth = 195 - 5/6

# Update contact names with titles
$ws.Range("C12").Value = "Jane Doe, Procurement Manager"
$ws.Range("C13").Value = "Ryan Smith, CTO"

# Capitalize "Not specified" -> "Not Specified"
$ws.Range("C19").Value = "Not Specified"
$ws.Range("C21").Value = "Not Specified"

# Expand renewal terms text
$ws.Range("C20").Value = "This Agreement shall automatically renew for successive one-year periods unless either party provides written notice of non-renewal at least 15 days prior to the expiration of the current term."
